$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update patient record values on row 2
$ws.Range("A2").Value = 3020
$ws.Range("E2").Value = 46200608020

# Leave selection on the edited cell, matching the saved workbook state
$ws.Activate()
$ws.Range("E2").Select()
